$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 201, shifting existing rows 201..258 down to 202..259
$ws.Rows.Item(201).Insert()

# Populate the newly inserted row 201 with the new record.
# Columns A,B,C,E,F,G,H,I,J,L,R keep the same values the old row 201 had
# (i.e. identical to what is now row 202), only D,K,M,N,O,P,Q,S,T differ.
$ws.Cells.Item(201, 1).Value = 7
$ws.Cells.Item(201, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(201, 3).Value = "Ñuble"
$ws.Cells.Item(201, 4).Value = 45120
$ws.Cells.Item(201, 5).Value = 16
$ws.Cells.Item(201, 6).Value = "Fruta"
$ws.Cells.Item(201, 7).Value = 100109
$ws.Cells.Item(201, 8).Value = "Uva"
$ws.Cells.Item(201, 9).Value = 100109001
$ws.Cells.Item(201, 10).Value = "Uva"
$ws.Cells.Item(201, 11).Value = "Crimpson Seedless"
$ws.Cells.Item(201, 12).Value = "Primera"
$ws.Cells.Item(201, 13).Value = 30
$ws.Cells.Item(201, 14).Value = 12000
$ws.Cells.Item(201, 15).Value = 12000
$ws.Cells.Item(201, 16).Value = 12000
$ws.Cells.Item(201, 17).Value = "$/bandeja 8 kilos"
$ws.Cells.Item(201, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(201, 19).Value = 1500
$ws.Cells.Item(201, 20).Value = 8
